$d = $word.ActiveDocument

# The document currently ends with a hyperlink paragraph (sciencedirect
# S1077201401902311) followed by one empty trailing paragraph, then the
# section properties. We need to add four more reference paragraphs,
# each containing a hyperlink, right before that trailing empty
# paragraph - the first of the four is also followed by a "." run.

$refParaIndex = 3
$refPara = $d.Paragraphs.Item($refParaIndex)

# Create a brand-new (bare) paragraph right after the reference paragraph;
# this pushes the existing empty trailing paragraph down by one and keeps
# it fully intact at the end of the document.
$refPara.Range.InsertParagraphAfter()

$url1 = "https://inspectorio.com/products/fabric-inspection-system#:~:text=Early%20detection%20of%20defects,more%20time%20and%20money%20saved"
$url2 = "https://kktmadhusanka.blogspot.com/2017/01/importance-of-fabric-inspection.html"
$url3 = "https://www.sciencedirect.com/science/article/abs/pii/S0030402614008523#preview-section-references"
$url4 = "https://www.tandfonline.com/doi/abs/10.1080/00405166908688985?journalCode=ttpr20"

# Insert the plain text for all four new paragraphs (plus the trailing
# "." after the first URL) in one shot, using paragraph marks to split
# them. Doing this before turning any of it into hyperlinks avoids the
# stray empty runs that Word's paragraph-insertion otherwise leaves
# behind.
$newPara = $d.Paragraphs.Item($refParaIndex + 1)
$insertPos = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$fullText = $url1 + "." + [char]13 + $url2 + [char]13 + $url3 + [char]13 + $url4
$insertPos.InsertAfter($fullText)

# Paragraph refParaIndex+1 now holds "url1.", refParaIndex+2 holds url2,
# refParaIndex+3 holds url3, refParaIndex+4 holds url4, and the original
# empty trailing paragraph follows right after as refParaIndex+5.

$p1 = $d.Paragraphs.Item($refParaIndex + 1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.Start + $url1.Length)
$d.Hyperlinks.Add($r1, $url1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url1) | Out-Null

$p2 = $d.Paragraphs.Item($refParaIndex + 2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$d.Hyperlinks.Add($r2, $url2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url2) | Out-Null

$p3 = $d.Paragraphs.Item($refParaIndex + 3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$d.Hyperlinks.Add($r3, $url3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url3) | Out-Null

$p4 = $d.Paragraphs.Item($refParaIndex + 4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$d.Hyperlinks.Add($r4, $url4, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url4) | Out-Null

Write-Output "Inserted 4 hyperlink paragraphs."
